$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.626.33'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.84%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.704.41'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '525.89'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.64'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.732.83'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.51%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.79'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +7.15%  '
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('E13').Value = '  +3.25%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.189.69'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.24%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '60.608.66'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.90%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.27'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.725.18'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.37%  '
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.50'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('E21').Value = '  +3.85%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.45'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.96%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.65%  '
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.996'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0₃0821'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.61%  '
$ws.Range('E29').Value = '  +2.67%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.82'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +9.12%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.04'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.83%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '149.72'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.26'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +7.00%  '
$ws.Range('E36').Value = '  +7.95%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.942'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.95%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.880'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.57%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.52'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +7.41%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '37.10'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.15%  '
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '281.19'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.612'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.77%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.996'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.143.88'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +7.69%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0987'
$ws.Range('D47').Style = "Normal"
$ws.Range('E48').Value = '  +6.09%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0538'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '10.53'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('E51').Value = '  +1.53%  '

Write-Output "Applied cryptos update"
